# Weekly price update: a new observation for "Pepino dulce" (week of
# 2023-05-04) is inserted at row 7 (right after the two "Terminal
# Hortofrutícola Agro Chillán" header-like rows that stay fixed at the
# top), pushing every existing record down by one row. The previously
# last record (old row 42) ends up at the new row 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; Excel shifts rows 7:42 down to
# 8:43 and extends the used range to A1:R43 automatically.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with this week's figures.
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 45050
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112043
$ws.Cells.Item(7, 7).Value = "Pepino dulce"
$ws.Cells.Item(7, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 50
$ws.Cells.Item(7, 11).Value = 14000
$ws.Cells.Item(7, 12).Value = 14000
$ws.Cells.Item(7, 13).Value = 14000
$ws.Cells.Item(7, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 778
$ws.Cells.Item(7, 17).Value = 18
$ws.Cells.Item(7, 18).Value = "Hortaliza"
